# Updates the "cryptos" price/volume table to the latest scraped snapshot.
# Cell values are plain display text (already formatted, e.g. "26.718.53" or
# "  -0.26%  "), so everything is written via .Value as a string. A handful of
# Price cells (column D) look like plain numbers (e.g. "218.03"); a leading
# apostrophe is prefixed -- Excel's normal force-text idiom -- so those stay
# text instead of being auto-coerced to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.718.53"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.633.46"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'218.03"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'18.96"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.861.14"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "1.621.27"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'64.05"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "26.690.63"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").Value = "'211.27"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("E23").Value = "  -5.04%  "
$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").Value = "'146.82"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "'6.99"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "'0.0500"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "1.263.06"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'0.0172"
$ws.Range("D38").Value = "'0.523"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'0.801"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").Value = "1.771.77"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "'59.63"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'0.0957"
$ws.Range("E51").Value = "  -0.73%  "
